$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ApplicationLogin")

# --- Row 2 ---
$ws.Range("A2").Value = 'Launch browser'
$ws.Range("B2").Value = 'startBrowser'
$ws.Range("C2").Value = 'NA'
$ws.Range("D2").Value = 'NA'
$ws.Range("E2").Value = 'NA'

# --- Row 3 ---
$ws.Range("A3").Value = 'Launch url in a browser'
$ws.Range("B3").Value = 'openApplication'
$ws.Range("C3").Value = 'NA'
$ws.Range("D3").Value = 'NA'
$ws.Range("E3").Value = 'NA'

# --- Rows 4-9 (Description/FunctionName/LocatorType/TestData first, LocatorValue after) ---
$ws.Range("A4").Value = 'wait for username'
$ws.Range("B4").Value = 'waitForElement'
$ws.Range("C4").Value = 'name'
$ws.Range("E4").Value = 10

$ws.Range("A5").Value = 'Enter username'
$ws.Range("B5").Value = 'typeAction'
$ws.Range("C5").Value = 'name'
$ws.Range("E5").Value = 'admin'

$ws.Range("A6").Value = 'wait for password'
$ws.Range("B6").Value = 'waitForElement'
$ws.Range("C6").Value = 'xpath'
$ws.Range("E6").Value = 10

$ws.Range("A7").Value = 'Enter password'
$ws.Range("B7").Value = 'typeAction'
$ws.Range("C7").Value = 'xpath'
$ws.Range("E7").Value = 'master'

$ws.Range("A8").Value = 'wait for login button'
$ws.Range("B8").Value = 'waitForElement'
$ws.Range("C8").Value = 'id'
$ws.Range("E8").Value = 10

$ws.Range("A9").Value = 'click login'
$ws.Range("B9").Value = 'clickAction'
$ws.Range("C9").Value = 'id'
$ws.Range("E9").Value = 'NA'

$ws.Range("D4").Value = 'username'
$ws.Range("D5").Value = 'username'
$ws.Range("D6").Value = '//input[@id=''password'']'
$ws.Range("D7").Value = '//input[@id=''password'']'
$ws.Range("D8").Value = 'btnsubmit'
$ws.Range("D9").Value = 'btnsubmit'

# --- Rows 10-13 ---
$ws.Range("A10").Value = 'wait for logout'
$ws.Range("B10").Value = 'waitForElement'
$ws.Range("C10").Value = 'xpath'
$ws.Range("D10").Value = '//a[@id=''logout'']'
$ws.Range("E10").Value = 10

$ws.Range("A11").Value = 'verify title'
$ws.Range("B11").Value = 'validateTitle'
$ws.Range("C11").Value = 'NA'
$ws.Range("D11").Value = 'NA'
$ws.Range("E11").Value = 'Dashboard « Stock Accounting'

$ws.Range("A12").Value = 'wait for logout'
$ws.Range("B12").Value = 'waitForElement'
$ws.Range("C12").Value = 'xpath'
$ws.Range("D12").Value = '//a[@id=''logout'']'
$ws.Range("E12").Value = 10

$ws.Range("A13").Value = 'click logout'
$ws.Range("B13").Value = 'clickAction'
$ws.Range("C13").Value = 'xpath'
$ws.Range("D13").Value = '//a[@id=''logout'']'
$ws.Range("E13").Value = 'NA'

# --- Rows 14-15 (Description/FunctionName/LocatorType/TestData first, LocatorValue after) ---
$ws.Range("A14").Value = 'wait for ok button'
$ws.Range("B14").Value = 'waitForElement'
$ws.Range("C14").Value = 'xpath'
$ws.Range("E14").Value = 10

$ws.Range("A15").Value = 'click ok button'
$ws.Range("B15").Value = 'clickAction'
$ws.Range("C15").Value = 'xpath'
$ws.Range("E15").Value = 'NA'

$ws.Range("D14").Value = '//button[normalize-space()=''OK!'']'
$ws.Range("D15").Value = '//button[normalize-space()=''OK!'']'

# --- Row 16 ---
$ws.Range("A16").Value = 'close browser'
$ws.Range("B16").Value = 'closeBrowser'
$ws.Range("C16").Value = 'NA'
$ws.Range("D16").Value = 'NA'
$ws.Range("E16").Value = 'NA'

# --- Font styling for the verified title value (E11) ---
$rngE11 = $ws.Range("E11")
$rngE11.Font.Name = "Consolas"
$rngE11.Font.Family = 3
$rngE11.Font.Size = 9
$rngE11.Font.Color = 2367776

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 23.25
$ws.Columns.Item(4).ColumnWidth = 31.25
$ws.Columns.Item(5).ColumnWidth = 31.25

# --- Page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("B3").Select()

Write-Host "Done"
